$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (preserve rich-text runs via Characters) ---
$ws.Range("A8").Characters(21, 2).Text = "16"
$ws.Range("C9").Characters(27, 8).Text = "4/15/2024"
$ws.Range("C9").Characters(47, 9).Text = "4/21/2024"

# --- Cells changing from numeric to text placeholder ("0" / "***.*") ---
$ws.Range("G14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("C18").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

# --- Cells changing from text placeholder to numeric ---
$ws.Range("D20").Value = 7
$ws.Range("C16").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = -85.714285714285
$ws.Range("E16").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("D29").Value = 2
$ws.Range("C16").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("E16").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("D30").Value = 2
$ws.Range("C16").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("E16").Copy()
$ws.Range("E30").PasteSpecial(-4122)

# --- Plain numeric value updates ---
$ws.Range("M14").Value = -20
$ws.Range("N15").Value = -50
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -42.857142857142
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -35.294117647058
$ws.Range("I16").Value = 49
$ws.Range("J16").Value = 54
$ws.Range("K16").Value = -9.259259259259
$ws.Range("L16").Value = -26.865671641791
$ws.Range("M16").Value = -32.876712328767
$ws.Range("N16").Value = -92.319749216300
$ws.Range("C17").Value = 9
$ws.Range("E17").Value = -10
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = -35
$ws.Range("I17").Value = 90
$ws.Range("J17").Value = 121
$ws.Range("K17").Value = -25.619834710743
$ws.Range("L17").Value = 2.272727272727
$ws.Range("M17").Value = 16.883116883116
$ws.Range("N17").Value = -72.560975609756
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -61.111111111111
$ws.Range("I18").Value = 40
$ws.Range("J18").Value = 43
$ws.Range("K18").Value = -6.976744186046
$ws.Range("L18").Value = -25.925925925925
$ws.Range("M18").Value = -50
$ws.Range("N18").Value = -91.543340380549
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -36.363636363636
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = -26.666666666666
$ws.Range("I19").Value = 111
$ws.Range("J19").Value = 128
$ws.Range("K19").Value = -13.28125
$ws.Range("L19").Value = -14.615384615384
$ws.Range("M19").Value = 52.054794520547
$ws.Range("N19").Value = -37.288135593220
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -81.818181818181
$ws.Range("J20").Value = 49
$ws.Range("K20").Value = -53.061224489795
$ws.Range("L20").Value = -34.285714285714
$ws.Range("M20").Value = -30.303030303030
$ws.Range("N20").Value = -91.417910447761
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = -50
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = -38.983050847457
$ws.Range("I21").Value = 326
$ws.Range("J21").Value = 405
$ws.Range("K21").Value = -19.506172839506
$ws.Range("L21").Value = -14.882506527415
$ws.Range("M21").Value = -5.780346820809
$ws.Range("N21").Value = -83.073727933541
$ws.Range("G22").Value = 2
$ws.Range("J22").Value = 8
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = -45.454545454545
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -60
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -27.272727272727
$ws.Range("I23").Value = 28
$ws.Range("J23").Value = 45
$ws.Range("K23").Value = -37.777777777777
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 86.666666666666
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 33.333333333333
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = -13.793103448275
$ws.Range("I24").Value = 214
$ws.Range("J24").Value = 220
$ws.Range("K24").Value = -2.727272727272
$ws.Range("L24").Value = -15.079365079365
$ws.Range("M24").Value = -10.084033613445
$ws.Range("C25").Value = 3
$ws.Range("E25").Value = 200
$ws.Range("I25").Value = 39
$ws.Range("J25").Value = 23
$ws.Range("K25").Value = 69.565217391304
$ws.Range("L25").Value = 62.5
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 18.181818181818
$ws.Range("F26").Value = 31
$ws.Range("H26").Value = -6.060606060606
$ws.Range("I26").Value = 129
$ws.Range("J26").Value = 179
$ws.Range("K26").Value = -27.932960893854
$ws.Range("L26").Value = -14.569536423841
$ws.Range("M26").Value = -40.825688073394
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = 66.666666666666
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -57.142857142857
$ws.Range("J28").Value = 22
$ws.Range("K28").Value = -36.363636363636
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 5
$ws.Range("K29").Value = -40
$ws.Range("M29").Value = -86.363636363636
$ws.Range("N29").Value = -96.153846153846
$ws.Range("G30").Value = 2
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = -40
$ws.Range("M30").Value = -82.352941176470
$ws.Range("N30").Value = -95.890410958904

$excel.CutCopyMode = 0
Write-Host "done"